$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "2022" column (S) to the table, mirroring the formatting of the
# existing 2021 column (R).

# Row 4 header: year label 2022, formatted like R4 (2021)
$ws.Range("R4").Copy() | Out-Null
$ws.Range("S4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("S4").Value = 2022

# Row 5 data: value 76.1, formatted like R5 (2021 value)
$ws.Range("R5").Copy() | Out-Null
$ws.Range("S5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("S5").Value = 76.1

$excel.CutCopyMode = $false

# Move the active selection, matching the saved workbook state
$ws.Range("P8").Select() | Out-Null
